$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarters: Dec-2018, Sep-2018)
$ws.Range("D1:E1").EntireColumn.Insert()

# Propagate number formatting from the (old D, now F) column into the two new columns
$ws.Range("F7:F35").Copy() | Out-Null
$ws.Range("D7:E35").PasteSpecial(-4122) | Out-Null
$ws.Range("F38:F66").Copy() | Out-Null
$ws.Range("D38:E66").PasteSpecial(-4122) | Out-Null
$ws.Range("F80:F102").Copy() | Out-Null
$ws.Range("D80:E102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the two newly inserted columns (D, E) with the new quarterly figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2770000
$ws.Range("E8").Value = 2349000
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 119000
$ws.Range("E14").Value = 100000
$ws.Range("D15").Value = 149000
$ws.Range("E15").Value = 138000
$ws.Range("D17").Value = 2273000
$ws.Range("E17").Value = 2090000
$ws.Range("D18").Value = 497000
$ws.Range("E18").Value = 259000
$ws.Range("D20").Value = -6000
$ws.Range("E20").Value = 4000
$ws.Range("D21").Value = 642000
$ws.Range("E21").Value = 403000
$ws.Range("D22").Value = 70000
$ws.Range("E22").Value = 69000
$ws.Range("D23").Value = 421000
$ws.Range("E23").Value = 194000
$ws.Range("D24").Value = 49000
$ws.Range("E24").Value = 39000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 372000
$ws.Range("E26").Value = 155000
$ws.Range("D27").Value = 364000
$ws.Range("E27").Value = 149000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -19000
$ws.Range("E29").Value = -2000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 6000
$ws.Range("E32").Value = -4000
$ws.Range("D33").Value = 345000
$ws.Range("E33").Value = 147000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 345000
$ws.Range("E35").Value = 147000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 656000
$ws.Range("E41").Value = 484000
$ws.Range("D42").Value = 4072000
$ws.Range("E42").Value = 4567000
$ws.Range("D43").Value = 9139000
$ws.Range("E43").Value = 7766000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 505000
$ws.Range("E45").Value = 531000
$ws.Range("D46").Value = 14372000
$ws.Range("E46").Value = 13348000
$ws.Range("D47").Value = 154000
$ws.Range("E47").Value = 127000
$ws.Range("D48").Value = 588000
$ws.Range("E48").Value = 594000
$ws.Range("D49").Value = 9320000
$ws.Range("E49").Value = 9542000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1988000
$ws.Range("E52").Value = 1991000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 26422000
$ws.Range("E54").Value = 25602000
$ws.Range("D57").Value = 1943000
$ws.Range("E57").Value = 1600000
$ws.Range("D58").Value = 251000
$ws.Range("E58").Value = 741000
$ws.Range("D59").Value = 11102000
$ws.Range("E59").Value = 10302000
$ws.Range("D60").Value = 13296000
$ws.Range("E60").Value = 12643000
$ws.Range("D61").Value = 6162000
$ws.Range("E61").Value = 5822000
$ws.Range("D62").Value = 2745000
$ws.Range("E62").Value = 2809000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 22271000
$ws.Range("E66").Value = 21340000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 2093000
$ws.Range("E72").Value = 2042000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 4151000
$ws.Range("E76").Value = 4262000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 345000
$ws.Range("E81").Value = 147000
$ws.Range("D83").Value = 151000
$ws.Range("E83").Value = 140000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 711000
$ws.Range("E89").Value = 562000
$ws.Range("D91").Value = -61000
$ws.Range("E91").Value = -68000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -53000
$ws.Range("E94").Value = -95000
$ws.Range("D96").Value = -97000
$ws.Range("E96").Value = -98000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -444000
$ws.Range("E100").Value = -457000
$ws.Range("D101").Value = -42000
$ws.Range("E101").Value = -13000
$ws.Range("D102").Value = 172000
$ws.Range("E102").Value = -3000

# Apply restated figures for previously-existing quarters (values revised in this update)
$ws.Range("F14").Value = 363000
$ws.Range("J14").Value = 535000
$ws.Range("F15").Value = 145000
$ws.Range("H15").Value = 178000
$ws.Range("I15").Value = 179000
$ws.Range("J15").Value = 177000
$ws.Range("H17").Value = 2322000
$ws.Range("I17").Value = 2084000
$ws.Range("H18").Value = 587000
$ws.Range("I18").Value = 256000
$ws.Range("H20").Value = -110000
$ws.Range("I20").Value = 14000
$ws.Range("H22").Value = 70000
$ws.Range("H24").Value = 44000
$ws.Range("H26").Value = 362000
$ws.Range("H27").Value = 355000
$ws.Range("H29").Value = -374000
$ws.Range("H32").Value = 110000
$ws.Range("I32").Value = -14000
